$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.939.94'
$ws.Range('E2').Value = '  -0.38%  '

# Row 3
$ws.Range('D3').Value = '3.117.83'
$ws.Range('E3').Value = '  +0.47%  '

# Row 4
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '578.30'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -0.68%  '

# Row 6
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '172.05'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +1.09%  '

# Row 7
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  +0.06%  '

# Row 8
$ws.Range('D8').Value = '3.114.07'
$ws.Range('E8').Value = '  +0.47%  '

# Row 9
$ws.Range('E9').Value = '  -0.87%  '

# Row 10
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '6.47'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -3.48%  '

# Row 11
$ws.Range('E11').Value = '  -1.91%  '

# Row 12
$ws.Range('E12').Value = '  +0.17%  '

# Row 13
$ws.Range('E13').Value = '  -2.03%  '

# Row 14
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '37.27'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  +0.72%  '

# Row 15
$ws.Range('E15').Value = '  -1.38%  '

# Row 16
$ws.Range('D16').Value = '3.637.08'
$ws.Range('E16').Value = '  +0.44%  '

# Row 17
$ws.Range('D17').Value = '66.897.09'
$ws.Range('E17').Value = '  -0.45%  '

# Row 18
$ws.Range('E18').Value = '  -1.47%  '

# Row 19
$ws.Range('D19').Value = '3.117.04'
$ws.Range('E19').Value = '  +0.34%  '

# Row 20
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '16.38'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +0.89%  '

# Row 21
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '476.78'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +0.89%  '

# Row 22
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '8.00'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +5.86%  '

# Row 23
$ws.Range('E23').Value = '  -0.59%  '

# Row 24
$ws.Range('E24').Value = '  +4.87%  '

# Row 25
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '84.06'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  +0.08%  '

# Row 26
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '2.31'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -2.63%  '

# Row 27
$ws.Range('E27').Value = '  -0.01%  '

# Row 28
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '10.01'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -2.69%  '

# Row 29
$ws.Range('B29').Value = 'NEARProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '7.89'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -2.86%  '

# Row 30
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '2.39'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  -1.63%  '

# Row 31
$ws.Range('E31').Value = '  -0.44%  '

# Row 32
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '28.61'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  +0.48%  '

# Row 33
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '0.116'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -0.17%  '

# Row 34
$ws.Range('D34').Value = '0.0₃0941'
$ws.Range('E34').Value = '  -7.72%  '

# Row 35
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -0.09%  '

# Row 36
$ws.Range('E36').Value = '  -1.20%  '

# Row 37
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '0.980'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  -3.51%  '

# Row 38
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '47.38'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  +0.09%  '

# Row 39
$ws.Range('E39').Value = '  -1.16%  '

# Row 40
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '50.02'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -0.99%  '

# Row 41
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '0.311'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -2.46%  '

# Row 42
$ws.Range('E42').Value = '  -1.76%  '

# Row 43
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '8.69'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -0.59%  '

# Row 44
$ws.Range('D44').Value = '2.812.16'
$ws.Range('E44').Value = '  +1.43%  '

# Row 45
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '0.0357'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  -2.38%  '

# Row 46
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '381.71'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  -3.86%  '

# Row 47
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '2.59'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -11.52%  '

# Row 48
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '136.17'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  +0.57%  '

# Row 49
$ws.Range('E49').Value = '  +0.05%  '

# Row 50
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '24.74'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -0.39%  '

# Row 51
$ws.Range('E51').Value = '  -2.44%  '
